$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# "Diagrama de Classe e Diagrama de Entidade e Relacionamento" (A9) is
# renamed to "Diagrama de Classe e Modelo Entidade-Relacionamento"
$ws.Range("A9").Value = "Diagrama de Classe e Modelo Entidade-Relacionamento"

# Match the author's final UI selection left in the saved file
[void]$ws.Range("D14").Select()
